# Updated symbol list on Tue Dec 20 22:13:07 UTC 2022 with GitHub Actions
#
# The "Price" (D), "Volume(1h)" (E) and "Hora" (G) columns are stored as
# TEXT in this sheet (e.g. "0.04030", "21"), even though most of them look
# numeric. Assigning a numeric-looking string straight to .Value would get
# re-interpreted as a Number and strip significant trailing zeros, so we
# force the target columns to Text number-format first (same as Excel's
# "Format Cells > Text") and then write the literal strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep columns D and G as Text so the literal strings (incl. trailing
# zeros like "0.1440") round-trip exactly instead of becoming numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# --- Column D: Price ---
$ws.Range("D2").Value = "251.13"
$ws.Range("D3").Value = "22.93"
$ws.Range("D4").Value = "5.421"
$ws.Range("D5").Value = "0.05669"
$ws.Range("D6").Value = "3.423"
$ws.Range("D7").Value = "6.378"
$ws.Range("D8").Value = "0.8139"
$ws.Range("D9").Value = "0.9324"
$ws.Range("D10").Value = "0.1440"
$ws.Range("D11").Value = "0.07447"
$ws.Range("D12").Value = "0.03110"
$ws.Range("D13").Value = "0.03066"
$ws.Range("D14").Value = "0.09355"
$ws.Range("D15").Value = "3.744"
$ws.Range("D16").Value = "0.001599"
$ws.Range("D17").Value = "0.04758"
$ws.Range("D18").Value = "0.0005793"
$ws.Range("D19").Value = "0.006405"
$ws.Range("D20").Value = "0.005043"
$ws.Range("D21").Value = "0.001032"
$ws.Range("D23").Value = "3.702"
$ws.Range("D24").Value = "2.182"
$ws.Range("D25").Value = "0.3303"
$ws.Range("D26").Value = "0.1312"
$ws.Range("D28").Value = "0.0003002"
$ws.Range("D40").Value = "0.04031"
$ws.Range("D41").Value = "0.006771"
$ws.Range("D42").Value = "0.1070"
$ws.Range("D44").Value = "0.008024"
$ws.Range("D45").Value = "0.00005806"
$ws.Range("D47").Value = "0.5003"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.01011"

# --- Column E: Volume(1h) label text tweaks ---
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

# --- Column G: Hora - every data row bumps from 21 to 22 ---
$ws.Range("G2:G51").Value = "22"
